$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 3 new rows for "reg-ack-template-part4" at rows 122-124,
# mirroring the existing "reg-ack-template-part3" block (rows 119-121).

$ws.Cells.Item(122, 1).Value = "reg-ack-template-part4"
$ws.Cells.Item(122, 2).Value = "Registration Acknowledgement Template - Part 4"
$ws.Cells.Item(122, 3).Value = "eng"
$ws.Cells.Item(122, 4).Value = $true
$ws.Cells.Item(122, 5).Value = "superadmin"
$ws.Cells.Item(122, 6).Value = "now()"

$ws.Cells.Item(123, 1).Value = "reg-ack-template-part4"
$ws.Cells.Item(123, 2).Value = "نموذج شكر التسجيل"
$ws.Cells.Item(123, 3).Value = "ara"
$ws.Cells.Item(123, 4).Value = $true
$ws.Cells.Item(123, 5).Value = "superadmin"
$ws.Cells.Item(123, 6).Value = "now()"

$ws.Cells.Item(124, 1).Value = "reg-ack-template-part4"
$ws.Cells.Item(124, 2).Value = "accusé de réception"
$ws.Cells.Item(124, 3).Value = "fra"
$ws.Cells.Item(124, 4).Value = $true
$ws.Cells.Item(124, 5).Value = "superadmin"
$ws.Cells.Item(124, 6).Value = "now()"

# Update the selection to mimic Excel's post-edit cursor position
$ws.Range("A125:XFD1048576").Select()
